$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Second run for Mtheo = 60, 0.85 : update the size estimate (was the
# 10x-too-small 0.1 um run) and relabel the Mtheo line to NA=60.
$ws.Range("B18").Value = "Taille estimé = 15.089146941064271 ± 0.6918618464405832 um et vrai Taille = 10 um"
$ws.Range("B32").Value = "Mtheo, NA = 60, 0.85"

# Leave the active selection where the author left it after pasting the
# image blob / adding the new run.
$ws.Range("C38").Select()
